$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update K column (최종점수 / final score) values
$ws.Range("K2").Value = 60.8
$ws.Range("K3").Value = 55
$ws.Range("K4").Value = 51.2
$ws.Range("K5").Value = 48.8

# Update N column (MACRO_SCORE) values
$ws.Range("N2").Value = 85.8724807945396
$ws.Range("N3").Value = 85.8724807945396
$ws.Range("N4").Value = 85.8724807945396
$ws.Range("N5").Value = 85.8724807945396
